# Update the "model 2 t/l ratio" calculation used throughout column AB.
# Old model: AB = 1.12 * AA^2 + 0.547 * AA + 0.066
# New model: AB = 1.1 * AA
# AB2 stands alone; AB3:AB66 and AB67:AB68 are shared-formula groups in the
# source sheet, so we mirror that grouping with ranged Formula assignment
# (one call per contiguous block) rather than one cell at a time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AB2").Formula = "=1.1 * AA2"
$ws.Range("AB3:AB66").Formula = "=1.1 * AA3"
$ws.Range("AB67:AB68").Formula = "=1.1 * AA67"

# Conditional formatting on column K: swap which rule (equals 0 / equals 1)
# carries which highlight color and put the "equals 1" (yellow) rule first/
# highest priority, matching the reorganized rule order.
$rng = $ws.Range("K1:K1048576")
$ruleZero = $rng.FormatConditions.Item(1)
$ruleOne = $rng.FormatConditions.Item(2)

$ruleZero.Formula1 = "=1"
$ruleOne.Formula1 = "=0"

$ruleZero.Font.Color = 393372
$ruleZero.Interior.Color = 13551615

$ruleOne.Font.Color = 22428
$ruleOne.Interior.Color = 10284031

$ruleZero.Priority = 1
$ruleOne.Priority = 2

# Leave the selection where the edited range is, as the author would after
# making this change.
$ws.Range("AB2:AB68").Select()
